$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text update (B1: jumlah_penetapan -> jumlah_penagihan) ---
$ws.Range("B1").Value = "jumlah_penagihan"

# --- New data rows ---
$ws.Range("A2").Value = "P.2.0021240.03.003"
$ws.Range("B2").Value = 120000
$ws.Range("A3").Value = "P.2.0021556.01.011"
$ws.Range("B3").Value = 2000
$ws.Range("A4").Value = "P.2.0001580.04.009."
$ws.Range("B4").Value = 400000

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 25.1640625
$ws.Columns.Item(2).ColumnWidth = 20.33203125

# --- Row heights for all 4 rows ---
$ws.Range("A1:B4").RowHeight = 18

# --- Fonts: build on scratch cells then paste formats onto target ranges to
#     avoid leaving intermediate styles applied to real cells ---
$ws.Range("Z1").Font.Name = "Arial"
$ws.Range("Z1").Font.Size = 14
$ws.Range("Z1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$ws.Range("Z2").Font.Name = "Helvetica Neue"
$ws.Range("Z2").Font.Size = 14
$ws.Range("Z2").Font.Color = 3355443
$ws.Range("Z2").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

$ws.Range("Z1:Z2").Clear()

# --- Selection ---
$ws.Range("B8").Select()
